$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix K11 (areaPoints for 冰川/Glacier, id 1007) - was a duplicate of K10, now unique
$ws.Range("K11").Value = "-4771|-40930|4348|-34123|-999|999"

# Update I14 (Loc for 帐篷/Tent, id 1010)
$ws.Range("I14").Value = "-6021.7|-50088|153.5"

# Delete row 15 (id 1011, 秘密基地/Secret Base) - all following rows shift up
$ws.Rows("15:15").Delete()

$ws.Range("K13").Select()
